# CCDD NTP Dosage form concepts.xlsx
# Add a new dose-form row: "oromucosal/laryngopharyngeal solution"
# Inserted as new row 91 (pushing the former rows 91-153 down to 92-154).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 91 (shifts everything below it down by one).
$ws.Rows.Item(91).Insert()

# Populate the new row's cells.
$ws.Range("B91").Value = "oromucosal/laryngopharyngeal solution"
$ws.Range("C91").Value = "Liquid preparation consisting of a solution intended for oromucosal or laryngopharyngeal use."
$ws.Range("E91").Value = "solution buccale/laryngopharyngée"

# Give E91 a plain (non-wrapped, General-format) style before tinting the
# font, so it doesn't inherit the neighbouring row's number format.
$ws.Range("E91").Style = "Normal"
$ws.Range("E91").Font.Color = 2236962   # RGB(0x22,0x22,0x22) -> FF222222

# Reflect the user's final selection/scroll position after the edit.
$ws.Range("E91").Select()

# Extend the hidden _FilterDatabase defined name by one row, as Excel does
# when a row is inserted inside/adjacent to the filtered range.
$n = $wb.Names.Item("_xlnm._FilterDatabase")
$n.RefersTo = "=Sheet1!`$E`$1:`$E`$155"
